$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
815948.8466242919,553903.4897101057,1160710.181512963
114717.946384911,85065.87905076986,151024.4192285407
148776.1234743129,106246.2524547077,197263.4535220726
154510.1013356931,118907.7896721767,205771.6466599592
390872.4106442083,258980.362916429,574741.2911992705
606880.8419751296,414314.4086768007,860624.4602578351
924344.6022658588,604470.7793365346,1263011.613799761
324275.2043504216,223957.0345758374,455600.9132062558
1057126.410802208,738840.6393578437,1420952.01392503
1035171.264536929,699267.8577558476,1446925.91437444
766586.6085540843,509706.8333459041,1090246.472679287
1371481.577681755,946325.6982386318,1837058.896405292
1165948.572484381,806300.827184155,1585026.115445837
205014.3129916234,148140.8884201075,285053.3513605037
389079.4361068338,259554.9399012513,563770.5013650062
762265.6030484788,498671.5026760732,1072088.966044433
148796.2469231778,106268.9699458452,197380.0781384635
1266437.489234939,844191.8897508249,1765092.177825624
509625.0798522626,334845.133003838,729401.0449461483
281437.8145121204,197298.9369640642,387534.2392097648
1198045.514401541,790119.1201548043,1664018.011755192
231750.535855732,158291.2947326502,329075.5248237757
392332.0440505766,262658.4911053977,567928.4198550846
266878.9867718245,182753.93987735,378829.0201037372
284895.2583476098,194827.8167130721,409251.7423977084
1615251.989634879,1065829.290860575,2240339.506126142
708529.7719326807,488376.6553981338,964515.48580973
261347.5460940035,181373.8186039565,364061.9060557179
661126.5335249294,445259.8214875035,925661.2496613846
1057787.297124199,729597.4347732206,1449998.530193612
997926.599332686,690952.1429656249,1390921.433959004
1163459.55117863,825084.0072865252,1575792.355092472
1134211.013164033,735090.4561967179,1616140.093486517
439992.1796270041,298680.5642920152,611515.8578083343
605183.0585923763,412842.6466459449,833966.4240149057
329303.7800887074,228208.0988412359,460976.724837868
1398996.562906179,945802.065668466,1963557.208743701
851745.9328587741,548332.9903191734,1208440.475573087
207416.0278052845,154550.0494706825,291067.7005485471
504652.5308022826,347112.3213318263,707383.3969430209
1143617.757840133,777632.4730633225,1602817.383408064
661180.1572990752,445339.0323934215,924840.590790672
665043.0275536403,462547.0965752712,929516.5778483752
1381445.565593534,943986.4401701103,1911865.780797239
917500.9548138493,646441.8308044226,1265247.096730388
1315891.055995775,894244.1706567031,1801595.498872063
766985.43717066,509746.6261386581,1095443.137453167
177274.4797560357,127447.4746467784,238921.9526097784
329113.9073687095,228147.3863543744,459605.6592446365
114114.5253083142,84424.58967766474,150569.380745057
785709.0024710825,535356.0541360932,1075299.267622101
838470.1877190369,543856.9746836675,1215270.458404636
201249.1380775836,145860.5039219384,280806.3515648408
282254.9273001692,196345.2786934428,406201.0346419477
1037267.124787318,706101.2751277878,1385097.30073972
997926.599332686,690952.1429656249,1390921.433959004
327645.0022719084,225806.8073685259,468721.9215611541
904333.3050928185,610271.7782148947,1259716.400569216
237473.340884903,163556.4629190774,339881.4507100929
1282941.499066801,861547.4266390188,1798421.359221304
1035301.735475383,698455.6868328867,1448030.987666231
1543061.96644725,1043564.072535655,2123786.919513411
237952.1292898582,164256.6169706851,336417.7384862658
869927.0548018506,578091.7028811517,1178681.729926297
1032620.930620987,680274.6940959555,1450758.799756194
918348.0463848049,617131.3385770199,1273612.297637243
248762.2230325001,168852.3233419045,362154.5311515725
1313927.536759081,918625.5164070942,1825688.992308316
207529.9648214925,150792.7989529555,285497.9836032013
341447.994752539,221079.1844486067,493284.8373962808
496711.2341149423,339488.3292196367,684236.0222262378
"@

$lines = $data -split "`n"
$startRow = 2

for ($i = 0; $i -lt $lines.Count; $i++) {
    $parts = $lines[$i].Trim().Split(",")
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = [double]$parts[0]
    $ws.Cells.Item($r, 2).Value2 = [double]$parts[1]
    $ws.Cells.Item($r, 3).Value2 = [double]$parts[2]
}
